$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: date in A6, error count in B6
$ws.Range("A6").Value = 45968
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = 45968
$ws.Range("B6").Value = 12

# Update selection to match the new active cell/selection state
$ws.Range("A6:B6").Select()
